$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (rows 28-34), matching columns:
# A=Item, B=Multiplicador, C=Pontos, D=Batalhas, E=ValorApostado, F=ValorAcumulado, G=Resultado
$rows = @(
    @("BonusPower", 0,                  0,     11,  40,  0,  "lose"),
    @("BonusPower", 2,                  720,   130, 50,  100, "win"),
    @("BonusPower", 0.13,               0,     56,  40,  0,  "lose"),
    @("SkipBoss",   0.7600000000000005, -619,  208, 20,  0,  "lose"),
    @("SkipBoss",   0.01,               10,    1,   200, 2,  "win"),
    @("SkipBoss",   2,                  1070,  107, 20,  40, "win"),
    @("BonusPower", 2,                  1070,  124, 21,  42, "win")
)

$startRow = 28
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
}
